# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "233.98") must be forced
# to text so Excel does not coerce the inline-string cell into a numeric cell.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D13",
    "D14",
    "D15",
    "D16",
    "D19",
    "D20",
    "D22",
    "D24",
    "D26",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50"
)

$updates = [ordered]@{
    "D2" = "37.774.20"
    "E2" = "  +1.08%  "
    "D3" = "2.082.71"
    "E3" = "  +0.72%  "
    "E4" = "  -0.02%  "
    "D5" = "233.98"
    "E5" = "  -0.38%  "
    "D6" = "0.624"
    "E6" = "  -0.18%  "
    "D7" = "58.79"
    "E7" = "  +2.74%  "
    "D9" = "0.390"
    "E9" = "  -1.04%  "
    "D10" = "0.0787"
    "E10" = "  +1.79%  "
    "E11" = "  +3.02%  "
    "D12" = "2.390.34"
    "E12" = "  +0.79%  "
    "D13" = "14.69"
    "E13" = "  +2.05%  "
    "D14" = "21.18"
    "E14" = "  +2.77%  "
    "D15" = "0.769"
    "E15" = "  -0.84%  "
    "D16" = "5.28"
    "E16" = "  +1.99%  "
    "D17" = "2.079.64"
    "E17" = "  +0.56%  "
    "D18" = "37.716.21"
    "E18" = "  +1.14%  "
    "D19" = "6.18"
    "E19" = "  +0.25%  "
    "D20" = "71.36"
    "E20" = "  +2.60%  "
    "D21" = "0.0₃0829"
    "E21" = "  +1.48%  "
    "D22" = "228.56"
    "E22" = "  +0.82%  "
    "E23" = "  -0.08%  "
    "D24" = "2.40"
    "E24" = "  -0.99%  "
    "E25" = "  -1.18%  "
    "D26" = "170.12"
    "E26" = "  +1.88%  "
    "E27" = "  +8.25%  "
    "D28" = "9.01"
    "E28" = "  +1.04%  "
    "D29" = "1.41"
    "E29" = "  +0.44%  "
    "D30" = "19.52"
    "E30" = "  +2.23%  "
    "E31" = "  +2.00%  "
    "D32" = "4.71"
    "E32" = "  +3.75%  "
    "D33" = "4.72"
    "E33" = "  +4.38%  "
    "D34" = "0.0628"
    "E34" = "  +2.10%  "
    "D35" = "2.51"
    "E35" = "  +1.81%  "
    "D36" = "3.44"
    "E36" = "  +3.64%  "
    "D37" = "1.83"
    "E37" = "  +2.37%  "
    "D38" = "1.00"
    "E38" = "  +0.04%  "
    "D39" = "5.39"
    "E39" = "  -3.72%  "
    "D40" = "0.0986"
    "E40" = "  +3.36%  "
    "E41" = "  +0.19%  "
    "D42" = "98.52"
    "E42" = "  +1.06%  "
    "D43" = "0.0214"
    "E43" = "  +1.10%  "
    "D44" = "1.456.89"
    "E44" = "  -1.98%  "
    "B45" = "FTXToken"
    "C45" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "D45" = "4.30"
    "E45" = "  +3.85%  "
    "B46" = "TrustWalletToken"
    "C46" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D46" = "1.17"
    "E46" = "  +0.76%  "
    "D47" = "16.01"
    "E47" = "  +6.42%  "
    "D48" = "1.07"
    "E48" = "  +4.36%  "
    "D49" = "7.39"
    "E49" = "  +2.90%  "
    "D50" = "3.03"
    "E50" = "  +2.60%  "
    "D51" = "2.275.00"
    "E51" = "  +0.67%  "
}

foreach ($addr in $updates.Keys) {
    $range = $ws.Range($addr)
    if ($textCells -contains $addr) {
        $range.NumberFormat = "@"
        $range.Value = $updates[$addr]
        $range.NumberFormat = "General"
    } else {
        $range.Value = $updates[$addr]
    }
}

Write-Host "Applied $($updates.Count) cell updates"
